# Updated symbol list on Sun Feb 12 05:31:35 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the crypto
# ranking table. Values are written with a leading apostrophe so Excel keeps
# them as text (matching the original inlineStr cell type) instead of
# reinterpreting "41.03" or "-0.06%" as numeric/percentage values. The Style
# reset afterwards clears the "quote prefix" formatting flag that Excel
# otherwise tags onto the cell, so no stray cell style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.06%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.53%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.246"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.27%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'0.44%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-0.19%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9183"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.80%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-2.54%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'13.78%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1831"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.19%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09200"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.01%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04258"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.15%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.05%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001263"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.75%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005750"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.17%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'3.354"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.312"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.18%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'1.22%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.404"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'12.90%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.48%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'2.93%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04069"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.90%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.29%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'7.51%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02471"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.32%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.00%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007847"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.80%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'0.98%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006821"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.98%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-1.86%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007706"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-7.78%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3052"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.50%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006723"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.83%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.22%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1699"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'592.28%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-2.38%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.22%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.22%"
$ws.Range("E51").Style = "Normal"
